# TireSearchAdvanced.xlsx - add a "Commercial Truck" positive-search scenario
# alongside the existing "Passenger car" one, and rename/re-layout a couple
# of header columns (old L/M "Remarks"/"Scenario Type" -> "Position"/"LoadRange",
# plus a new trailing "Remarks" column N).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header row (row 1) ----------------------------------------------
$ws.Range("L1").Value = "Position"
$ws.Range("M1").Value = "LoadRange"
$ws.Range("N1").Value = "Remarks"

# ---- Row 2 (existing Passenger-car "Yes" row) -------------------------
$ws.Range("G2").Value = "BF Goodrich"
$ws.Range("I2").Value = "18"
$ws.Range("K2").Value = "T"
$ws.Range("L2").ClearContents()
$ws.Range("M2").ClearContents()
$ws.Range("N2").Value = "Passenger car - Positive search"

# ---- Row 3 (existing Passenger-car "No" row) ---------------------------
$ws.Range("L3").ClearContents()
$ws.Range("M3").ClearContents()
$ws.Range("N3").Value = "Passenger car - Positive search"

# ---- Row 5 (new Commercial-Truck row, fully populated) -----------------
$ws.Range("A5").Value = "No"
$ws.Range("B5").Value = "CT"
$ws.Range("C5").Value = 1002750
$ws.Range("D5").Value = 1276063
$ws.Range("E5").Value = "275"
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = "Michelin"
$ws.Range("L5").Value = "Drive"
$ws.Range("M5").Value = "G"
$ws.Range("N5").Value = "Commercial Truck - Positive Search"

# ---- Row 6 (new Commercial-Truck row, partially populated) -------------
$ws.Range("A6").Value = "No"
$ws.Range("B6").Value = "CT"
$ws.Range("C6").Value = 1002750
$ws.Range("D6").Value = 1276063
$ws.Range("N6").Value = "Commercial Truck - Positive Search"

# ---- Column widths -------------------------------------------------
# Columns A/B got narrower, the old "Remarks" column (12) reverts to the
# default width, and the vacated/ new columns 13 (LoadRange) & 14 (Remarks)
# pick up the widths the old columns 12/13 used to have.
$ws.Columns.Item(1).ColumnWidth = 3.67
$ws.Columns.Item(2).ColumnWidth = 4.5
$ws.Columns.Item(13).ColumnWidth = 10.17
$ws.Columns.Item(14).ColumnWidth = 32

# ---- Selection ----------------------------------------------------
$ws.Range("N5").Select()
